$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.993.15"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "2.653.95"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.42%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "3.113.52"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "59.007.04"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "2.642.70"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "339.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("E19").Value = "  -4.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.419"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").Value = "0.0₃0803"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.892"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.869"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.616"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0968"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "2.035.24"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
